$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-14
$ws.Range("C2:C14").Value = 46070

# Rows 10-13 get new values for A (Beteckning), B (Datum), G (Area)
# The new case "A 3402-2026" moved to the top (row 10), shifting the
# previous rows 10-12 down to rows 11-13, and the former row 13 moved out.

$ws.Range("A10").Value = "A 3402-2026"
$ws.Range("B10").Value = 46042.39047453704
$ws.Range("G10").Value = 5.5

$ws.Range("A11").Value = "A 25015-2023"
$ws.Range("B11").Value = 45085.6989699074
$ws.Range("G11").Value = 1.8

$ws.Range("A12").Value = "A 19922-2025"
$ws.Range("B12").Value = 45771.63034722222
$ws.Range("G12").Value = 10.1

$ws.Range("A13").Value = "A 60024-2025"
$ws.Range("B13").Value = 45992
$ws.Range("G13").Value = 1.1
